# Ajuste no algoritmo para considerar o preço de fechamento
# Atualiza os resultados da simulação (planilha SimulacaoPeloLote)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colunas H, I (decimais "D.Inicial"/"D.Final") e P (datas "Data Melhor Saldo")
# guardam texto que se parece com numero/data. Forcamos formato de texto
# antes de atribuir o valor para preservar o tipo de celula original (texto).
# (Ranges multi-area nao aplicam a todas as areas, entao setamos celula a celula.)
$ws.Range("H2:H12").NumberFormat = "@"
foreach ($addr in @("I2","I3","I4","I5","I6","I9","I11","I12")) {
    $ws.Range($addr).NumberFormat = "@"
}
foreach ($addr in @("P2","P3","P5","P7","P8","P12")) {
    $ws.Range($addr).NumberFormat = "@"
}

# Linha 2 - PETRE301
$ws.Range("H2").Value = "0.9992"
$ws.Range("I2").Value = "0.5883"
$ws.Range("N2").Value = "R$ -710.32"
$ws.Range("O2").Value = "R$ -17.93"
$ws.Range("P2").Value = "2025-04-04"

# Linha 3 - PETRE306
$ws.Range("H3").Value = "0.9976"
$ws.Range("I3").Value = "0.4645"
$ws.Range("N3").Value = "R$ -811.25"
$ws.Range("O3").Value = "R$ -99.83"
$ws.Range("P3").Value = "2025-04-02"

# Linha 4 - PETRE312
$ws.Range("H4").Value = "0.9952"
$ws.Range("I4").Value = "0.3847"
$ws.Range("K4").Value = 200
$ws.Range("M4").Value = 2
$ws.Range("N4").Value = "R$ -958.81"
$ws.Range("O4").Value = "R$ -77.41"

# Linha 5 - PETRE316
$ws.Range("H5").Value = "0.9911"
$ws.Range("I5").Value = "0.3107"
$ws.Range("J5").Value = "ITM → ATM"
$ws.Range("N5").Value = "R$ -595.32"
$ws.Range("O5").Value = "R$ -69.38"
$ws.Range("P5").Value = "2025-04-02"

# Linha 6 - PETRE321
$ws.Range("H6").Value = "0.9843"
$ws.Range("I6").Value = "0.2447"
$ws.Range("K6").Value = 200
$ws.Range("M6").Value = 3
$ws.Range("N6").Value = "R$ -640.67"
$ws.Range("O6").Value = "R$ -108.90"

# Linha 7 - PETRF303
$ws.Range("H7").Value = "0.6628"
$ws.Range("N7").Value = "R$ 1050.08"
$ws.Range("O7").Value = "R$ 1050.08"
$ws.Range("P7").Value = "2025-06-18"

# Linha 8 - PETRF321
$ws.Range("H8").Value = "0.4271"
$ws.Range("N8").Value = "R$ 241.96"
$ws.Range("O8").Value = "R$ 770.78"
$ws.Range("P8").Value = "2025-06-10"

# Linha 9 - PETRF326
$ws.Range("H9").Value = "0.3690"
$ws.Range("I9").Value = "0.9982"
$ws.Range("L9").Value = 60
$ws.Range("N9").Value = "R$ 308.47"
$ws.Range("O9").Value = "R$ 968.36"

# Linha 10 - PETRF331
$ws.Range("H10").Value = "0.5522"
$ws.Range("L10").Value = 252
$ws.Range("N10").Value = "R$ 2047.65"
$ws.Range("O10").Value = "R$ 2047.65"

# Linha 11 - PETRF342
$ws.Range("H11").Value = "0.2198"
$ws.Range("I11").Value = "0.6801"
$ws.Range("N11").Value = "R$ 570.84"
$ws.Range("O11").Value = "R$ 634.14"

# Linha 12 - PETRF376
$ws.Range("H12").Value = "0.2016"
$ws.Range("I12").Value = "0.1700"
$ws.Range("N12").Value = "R$ 713.64"
$ws.Range("O12").Value = "R$ 713.64"
$ws.Range("P12").Value = "2025-06-18"
